# Apply updated probability values to the "team specific matrix" sheet.
# The workbook tracks per-state transition probabilities; this edit refreshes
# the simulated values after adding more games / speeding up simulate-game logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1836065573770492
$ws.Range("C2").Value = 0.5704918032786885
$ws.Range("J2").Value = 0.006557377049180328
$ws.Range("P2").Value = 0.1475409836065574
$ws.Range("S2").Value = 0.09180327868852459
$ws.Range("C3").Value = 0.005681818181818182
$ws.Range("J3").Value = 0.02840909090909091
$ws.Range("P3").Value = 0.6988636363636364
$ws.Range("S3").Value = 0.2670454545454545
$ws.Range("J4").Value = 0.08823529411764706
$ws.Range("P4").Value = 0.6764705882352942
$ws.Range("S4").Value = 0.2352941176470588
$ws.Range("B6").Value = 0.06097560975609756
$ws.Range("D6").Value = 0.01626016260162602
$ws.Range("E6").Value = 0.004065040650406504
$ws.Range("F6").Value = 0.06097560975609756
$ws.Range("J6").Value = 0.3130081300813008
$ws.Range("O6").Value = 0.02032520325203252
$ws.Range("Q6").Value = 0.1097560975609756
$ws.Range("R6").Value = 0.06910569105691057
$ws.Range("S6").Value = 0.3455284552845528
$ws.Range("B7").Value = 0.09826589595375723
$ws.Range("D7").Value = 0.0115606936416185
$ws.Range("F7").Value = 0.05202312138728324
$ws.Range("J7").Value = 0.138728323699422
$ws.Range("O7").Value = 0.02312138728323699
$ws.Range("Q7").Value = 0.1445086705202312
$ws.Range("R7").Value = 0.09248554913294797
$ws.Range("S7").Value = 0.4393063583815029
$ws.Range("B8").Value = 0.1042183622828784
$ws.Range("D8").Value = 0.01240694789081886
$ws.Range("F8").Value = 0.0620347394540943
$ws.Range("J8").Value = 0.09181141439205956
$ws.Range("O8").Value = 0.03225806451612903
$ws.Range("Q8").Value = 0.1588089330024814
$ws.Range("R8").Value = 0.08436724565756824
$ws.Range("S8").Value = 0.4540942928039702
$ws.Range("B9").Value = 0.09743589743589744
$ws.Range("D9").Value = 0.005128205128205128
$ws.Range("F9").Value = 0.09230769230769231
$ws.Range("J9").Value = 0.1128205128205128
$ws.Range("O9").Value = 0.01025641025641026
$ws.Range("Q9").Value = 0.158974358974359
$ws.Range("R9").Value = 0.07179487179487179
$ws.Range("S9").Value = 0.4512820512820513
$ws.Range("B10").Value = 0.1377459749552773
$ws.Range("D10").Value = 0.02057245080500894
$ws.Range("E10").Value = 0.001788908765652952
$ws.Range("F10").Value = 0.08228980322003578
$ws.Range("J10").Value = 0.09838998211091235
$ws.Range("O10").Value = 0.01699463327370304
$ws.Range("Q10").Value = 0.165474060822898
$ws.Range("R10").Value = 0.08407871198568873
$ws.Range("S10").Value = 0.3926654740608229
$ws.Range("G11").Value = 0.1745454545454546
$ws.Range("J11").Value = 0.07272727272727272
$ws.Range("K11").Value = 0.2327272727272727
$ws.Range("L11").Value = 0.509090909090909
$ws.Range("S11").Value = 0.01090909090909091
$ws.Range("G12").Value = 0.7482014388489209
$ws.Range("J12").Value = 0.1870503597122302
$ws.Range("K12").Value = 0.01438848920863309
$ws.Range("L12").Value = 0.01438848920863309
$ws.Range("S12").Value = 0.03597122302158273
$ws.Range("G13").Value = 0.6279069767441861
$ws.Range("J13").Value = 0.3255813953488372
$ws.Range("S13").Value = 0.04651162790697674
$ws.Range("F15").Value = 0.02392344497607655
$ws.Range("H15").Value = 0.1339712918660287
$ws.Range("I15").Value = 0.09090909090909091
$ws.Range("J15").Value = 0.3397129186602871
$ws.Range("K15").Value = 0.03827751196172249
$ws.Range("M15").Value = 0.01913875598086124
$ws.Range("O15").Value = 0.09090909090909091
$ws.Range("S15").Value = 0.2631578947368421
$ws.Range("F16").Value = 0.01587301587301587
$ws.Range("H16").Value = 0.1904761904761905
$ws.Range("I16").Value = 0.09523809523809523
$ws.Range("J16").Value = 0.3333333333333333
$ws.Range("K16").Value = 0.126984126984127
$ws.Range("M16").Value = 0.03174603174603174
$ws.Range("O16").Value = 0.08465608465608465
$ws.Range("S16").Value = 0.1216931216931217
$ws.Range("F17").Value = 0.02694610778443114
$ws.Range("H17").Value = 0.2005988023952096
$ws.Range("I17").Value = 0.09281437125748503
$ws.Range("J17").Value = 0.437125748502994
$ws.Range("K17").Value = 0.07784431137724551
$ws.Range("M17").Value = 0.02395209580838323
$ws.Range("O17").Value = 0.06287425149700598
$ws.Range("S17").Value = 0.07784431137724551
$ws.Range("F18").Value = 0.03428571428571429
$ws.Range("H18").Value = 0.16
$ws.Range("I18").Value = 0.08571428571428572
$ws.Range("J18").Value = 0.4685714285714286
$ws.Range("K18").Value = 0.08
$ws.Range("M18").Value = 0.01714285714285714
$ws.Range("O18").Value = 0.07428571428571429
$ws.Range("S18").Value = 0.08
$ws.Range("F19").Value = 0.03456998313659359
$ws.Range("H19").Value = 0.2141652613827993
$ws.Range("I19").Value = 0.09359190556492411
$ws.Range("J19").Value = 0.3659359190556493
$ws.Range("K19").Value = 0.1112984822934233
$ws.Range("M19").Value = 0.02107925801011805
$ws.Range("O19").Value = 0.05986509274873524
$ws.Range("S19").Value = 0.09949409780775717

